$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DATE column (B) for every data row (2-80) moves from 2024-10-01 to
# 2024-10-05 now that the prediction/training run has completed. Dates are
# stored as plain text in this sheet, so a leading apostrophe is used to
# force a literal text entry (preventing Excel from auto-converting the
# date-shaped string into a real date serial number).
for ($row = 2; $row -le 80; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Value2 -eq "2024-10-01") {
        $cell.Value2 = "'2024-10-05"
    }
}
